$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 9702.25
$ws.Range("I43").Value = 8999.944
$ws.Range("J43").Value = 10404.556
$ws.Range("K43").Value = 8999.944
$ws.Range("L43").Value = 10404.556
$ws.Range("M43").Value = -8930.944
$ws.Range("N43").Value = -10542.556
$ws.Range("H70").Value = 4008.75
$ws.Range("I70").Value = 1715.8334
$ws.Range("J70").Value = 6301.6665
$ws.Range("K70").Value = 5147.5002
$ws.Range("L70").Value = 18904.9995
$ws.Range("M70").Value = -4877.5002
$ws.Range("N70").Value = -19444.9995
$ws.Range("H73").Value = 4008.75
$ws.Range("I73").Value = 1715.8334
$ws.Range("J73").Value = 6301.6665
$ws.Range("K73").Value = 5147.5002
$ws.Range("L73").Value = 18904.9995
$ws.Range("M73").Value = -4211.5002
$ws.Range("N73").Value = -20776.9995
$ws.Range("H88").Value = 6128.222
$ws.Range("I88").Value = 5500
$ws.Range("K88").Value = 5500
$ws.Range("M88").Value = -5094
$ws.Range("H91").Value = 6128.222
$ws.Range("I91").Value = 5500
$ws.Range("K91").Value = 5500
$ws.Range("M91").Value = -4096
$ws.Range("H132").Value = 1196.4755
$ws.Range("I132").Value = 1183.5
$ws.Range("J132").Value = 1296.5714
$ws.Range("K132").Value = 3550.5
$ws.Range("L132").Value = 3889.7142
$ws.Range("M132").Value = -1020.5
$ws.Range("N132").Value = -8949.7142
$ws.Range("H137").Value = 41670510
$ws.Range("I137").Value = 125003896
$ws.Range("K137").Value = 375011688
$ws.Range("M137").Value = -375009138
$ws.Range("H138").Value = 3263.6765
$ws.Range("J138").Value = 3196.4443
$ws.Range("L138").Value = 9589.332900000001
$ws.Range("N138").Value = -19869.3329

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18520844
$ws.Range("I32").Value = 18869746
$ws.Range("K32").Value = 18869746
$ws.Range("M32").Value = -18869459
$ws.Range("H74").Value = 12822092
$ws.Range("J74").Value = 2400
$ws.Range("L74").Value = 2400
$ws.Range("H77").Value = 12822092
$ws.Range("J77").Value = 2400
$ws.Range("L77").Value = 12000
$ws.Range("H110").Value = 9670.833000000001
$ws.Range("I110").Value = 3999.6667
$ws.Range("K110").Value = 3999.6667
$ws.Range("M110").Value = -1954.6667
$ws.Range("H122").Value = 47620184
$ws.Range("I122").Value = 1152.6316
$ws.Range("K122").Value = 3457.8948
$ws.Range("M122").Value = -1007.8948
$ws.Range("H132").Value = 3722.597
$ws.Range("I132").Value = 1481.6545
$ws.Range("K132").Value = 4444.9635
$ws.Range("M132").Value = -1914.9635
$ws.Range("N74").Value = -4148
$ws.Range("N77").Value = -20736

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 10000
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("H134").Value = 3461.158
$ws.Range("I134").Value = 1394.9333
$ws.Range("K134").Value = 4184.7999
$ws.Range("M134").Value = -1649.7999
$ws.Range("N46").Value = -10596

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2421.9412
$ws.Range("I16").Value = 1329.3
$ws.Range("J16").Value = 3982.8572
$ws.Range("K16").Value = 1329.3
$ws.Range("L16").Value = 3982.8572
$ws.Range("M16").Value = -1042.3
$ws.Range("N16").Value = -4556.8572
$ws.Range("H22").Value = 3026.4285
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3026.4285
$ws.Range("K22").Value = 0
$ws.Range("N22").Value = -3726.4285
$ws.Range("H31").Value = 34766.305
$ws.Range("I31").Value = 2215.182
$ws.Range("J31").Value = 99868.55
$ws.Range("K31").Value = 2215.182
$ws.Range("L31").Value = 99868.55
$ws.Range("M31").Value = -1920.182
$ws.Range("N31").Value = -100458.55
$ws.Range("H34").Value = 34766.305
$ws.Range("I34").Value = 2215.182
$ws.Range("J34").Value = 99868.55
$ws.Range("K34").Value = 2215.182
$ws.Range("L34").Value = 99868.55
$ws.Range("M34").Value = -2013.182
$ws.Range("N34").Value = -100272.55
$ws.Range("H47").Value = 28999
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 28999
$ws.Range("K47").Value = 0
$ws.Range("N47").Value = -30131
$ws.Range("H86").Value = 8844.272000000001
$ws.Range("I86").Value = 3649.5
$ws.Range("J86").Value = 9998.666999999999
$ws.Range("K86").Value = 3649.5
$ws.Range("L86").Value = 9998.666999999999
$ws.Range("M86").Value = -2526.5
$ws.Range("N86").Value = -12244.667
$ws.Range("H89").Value = 8844.272000000001
$ws.Range("I89").Value = 3649.5
$ws.Range("J89").Value = 9998.666999999999
$ws.Range("K89").Value = 18247.5
$ws.Range("L89").Value = 49993.335
$ws.Range("M89").Value = -12631.5
$ws.Range("N89").Value = -61225.335
$ws.Range("H113").Value = 2421.9412
$ws.Range("I113").Value = 1329.3
$ws.Range("J113").Value = 3982.8572
$ws.Range("K113").Value = 1329.3
$ws.Range("L113").Value = 3982.8572
$ws.Range("M113").Value = 840.7
$ws.Range("N113").Value = -8322.8572
$ws.Range("H132").Value = 2785.7778
$ws.Range("I132").Value = 613.6
$ws.Range("K132").Value = 1840.8
$ws.Range("M132").Value = 689.1999999999998
$ws.Range("L22").Value = 3026.4285
$ws.Range("L47").Value = 28999
$ws.Range("M22").ClearContents()
$ws.Range("M47").ClearContents()

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 58
$ws.Range("I13").Value = 55
$ws.Range("J13").Value = 59.5
$ws.Range("K13").Value = 165
$ws.Range("L13").Value = 178.5
$ws.Range("M13").Value = 3
$ws.Range("H50").Value = 41676200
$ws.Range("I50").Value = 66668836
$ws.Range("K50").Value = 200006508
$ws.Range("M50").Value = -200006027
$ws.Range("H53").Value = 41676200
$ws.Range("I53").Value = 66668836
$ws.Range("K53").Value = 200006508
$ws.Range("M53").Value = -200006027
$ws.Range("N13").Value = -514.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1010
$ws.Range("I113").Value = 1022.1111
$ws.Range("J113").Value = 901
$ws.Range("K113").Value = 1022.1111
$ws.Range("L113").Value = 901
$ws.Range("M113").Value = 1147.8889
$ws.Range("N113").Value = -5241
$ws.Range("H122").Value = 4929
$ws.Range("I122").Value = 1700
$ws.Range("J122").Value = 8158
$ws.Range("K122").Value = 5100
$ws.Range("L122").Value = 24474
$ws.Range("M122").Value = -2650
$ws.Range("N122").Value = -29374
$ws.Range("H132").Value = 876603.4399999999
$ws.Range("I132").Value = 1192999.4
$ws.Range("K132").Value = 3578998.2
$ws.Range("M132").Value = -3576468.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6188.1113
$ws.Range("I7").Value = 3269.4119
$ws.Range("K7").Value = 3269.4119
$ws.Range("M7").Value = -3157.4119
$ws.Range("H16").Value = 2875.2
$ws.Range("I16").Value = 2300.125
$ws.Range("J16").Value = 5175.5
$ws.Range("K16").Value = 2300.125
$ws.Range("L16").Value = 5175.5
$ws.Range("M16").Value = -2130.125
$ws.Range("N16").Value = -5515.5
$ws.Range("H22").Value = 6693.95
$ws.Range("I22").Value = 3224.5
$ws.Range("J22").Value = 9006.916999999999
$ws.Range("K22").Value = 3224.5
$ws.Range("L22").Value = 9006.916999999999
$ws.Range("M22").Value = -2929.5
$ws.Range("N22").Value = -9596.916999999999
$ws.Range("H27").Value = 6693.95
$ws.Range("I27").Value = 3224.5
$ws.Range("J27").Value = 9006.916999999999
$ws.Range("K27").Value = 3224.5
$ws.Range("L27").Value = 9006.916999999999
$ws.Range("M27").Value = -3117.5
$ws.Range("N27").Value = -9220.916999999999
$ws.Range("H40").Value = 13652.429
$ws.Range("I40").Value = 15012.4
$ws.Range("K40").Value = 15012.4
$ws.Range("M40").Value = -14876.4
$ws.Range("H46").Value = 2736.5334
$ws.Range("I46").Value = 667
$ws.Range("J46").Value = 3253.9167
$ws.Range("K46").Value = 667
$ws.Range("L46").Value = 3253.9167
$ws.Range("M46").Value = -479
$ws.Range("N46").Value = -3629.9167
$ws.Range("H61").Value = 3264.7222
$ws.Range("I61").Value = 1297
$ws.Range("J61").Value = 13103.333
$ws.Range("K61").Value = 1297
$ws.Range("L61").Value = 13103.333
$ws.Range("M61").Value = -1095
$ws.Range("N61").Value = -13507.333
$ws.Range("H68").Value = 3805.5625
$ws.Range("I68").Value = 2917.6296
$ws.Range("J68").Value = 8600.4
$ws.Range("K68").Value = 2917.6296
$ws.Range("L68").Value = 8600.4
$ws.Range("M68").Value = -2168.6296
$ws.Range("N68").Value = -10098.4
$ws.Range("H71").Value = 3805.5625
$ws.Range("I71").Value = 2917.6296
$ws.Range("J71").Value = 8600.4
$ws.Range("K71").Value = 14588.148
$ws.Range("L71").Value = 43002
$ws.Range("M71").Value = -10844.148
$ws.Range("N71").Value = -50490
$ws.Range("H111").Value = 73684
$ws.Range("J111").Value = 73684
$ws.Range("L111").Value = 73684
$ws.Range("H113").Value = 3264.7222
$ws.Range("I113").Value = 1297
$ws.Range("J113").Value = 13103.333
$ws.Range("K113").Value = 1297
$ws.Range("L113").Value = 13103.333
$ws.Range("M113").Value = 873
$ws.Range("N113").Value = -17443.333
$ws.Range("H122").Value = 7117.3335
$ws.Range("I122").Value = 3424.75
$ws.Range("K122").Value = 10274.25
$ws.Range("M122").Value = -7824.25
$ws.Range("H126").Value = 6188.1113
$ws.Range("I126").Value = 3269.4119
$ws.Range("K126").Value = 9808.235700000001
$ws.Range("M126").Value = -7338.235700000001
$ws.Range("H132").Value = 8457.130999999999
$ws.Range("I132").Value = 4136.727
$ws.Range("J132").Value = 12417.5
$ws.Range("K132").Value = 12410.181
$ws.Range("L132").Value = 37252.5
$ws.Range("M132").Value = -9880.181
$ws.Range("N132").Value = -42312.5
$ws.Range("H136").Value = 10961.467
$ws.Range("I136").Value = 4578.3335
$ws.Range("J136").Value = 15216.889
$ws.Range("K136").Value = 13735.0005
$ws.Range("L136").Value = 45650.667
$ws.Range("M136").Value = -11185.0005
$ws.Range("N136").Value = -50750.667
$ws.Range("N111").Value = -81864

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 615.9167
$ws.Range("I113").Value = 253
$ws.Range("J113").Value = 736.8889
$ws.Range("K113").Value = 759
$ws.Range("L113").Value = 2210.6667
$ws.Range("M113").Value = 1411
$ws.Range("N113").Value = -6550.6667
$ws.Range("H136").Value = 1732.1143
$ws.Range("I136").Value = 1224.0883
$ws.Range("K136").Value = 3672.2649
$ws.Range("M136").Value = -1122.2649
